# Apply the SKU correction: "D184T22BKD03" -> "D184T22BKD04"
# and update the active selection to B39 (as recorded in the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product SKU text in cell A38.
$ws.Range("A38").Value = "D184T22BKD04"

# Update the selected cell to reflect the new active cell (B39) saved with the sheet.
$ws.Range("B39").Select()
